$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.044.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "'1.646.01"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.12%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'206.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.07%  "
$ws.Range("D6").Value = "'0.5191"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.71%  "
$ws.Range("D7").Value = "'1.003"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'0.2574"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.06%  "
$ws.Range("D9").Value = "'0.06248"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.79%  "
$ws.Range("D10").Value = "'20.73"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.46%  "
$ws.Range("D11").Value = "'0.07565"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.83%  "
$ws.Range("D12").Value = "'1.640.14"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.53%  "
$ws.Range("D13").Value = "'4.382"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.47%  "
$ws.Range("D14").Value = "'0.5391"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.08%  "
$ws.Range("D15").Value = "'66.07"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.02%  "
$ws.Range("D16").Value = "'0.0₅7928"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.38%  "
$ws.Range("D17").Value = "'26.009.24"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.29%  "
$ws.Range("D18").Value = "'1.005"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.17%  "
$ws.Range("D19").Value = "'4.669"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.05%  "
$ws.Range("D20").Value = "'187.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.98%  "
$ws.Range("D21").Value = "'10.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.49%  "
$ws.Range("D22").Value = "'6.126"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").Value = "'1.004"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").Value = "'148.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.01%  "
$ws.Range("D25").Value = "'0.1215"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.72%  "
$ws.Range("D26").Value = "'7.350"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.12%  "
$ws.Range("D27").Value = "'15.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.54%  "
$ws.Range("E28").Value = "  +2.46%  "
$ws.Range("D29").Value = "'0.06012"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.24%  "
$ws.Range("D30").Value = "'1.241"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.23%  "
$ws.Range("D31").Value = "'3.440"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.72%  "
$ws.Range("D32").Value = "'3.396"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.29%  "
$ws.Range("D33").Value = "'1.627"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.14%  "
$ws.Range("D34").Value = "'0.9771"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.55%  "
$ws.Range("D35").Value = "'2.384"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.99%  "
$ws.Range("D36").Value = "'2.730"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.08%  "
$ws.Range("D37").Value = "'0.5855"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.52%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.01591"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.89%  "
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "'1.085.22"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.91%  "
$ws.Range("D40").Value = "'5.976"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.80%  "
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("D42").Value = "'0.8433"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.45%  "
$ws.Range("D43").Value = "'100.30"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.49%  "
$ws.Range("D44").Value = "'1.798.82"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.77%  "
$ws.Range("D45").Value = "'0.0₈105"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.13%  "
$ws.Range("D46").Value = "'1.005"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.18%  "
$ws.Range("D47").Value = "'54.81"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.85%  "
$ws.Range("D48").Value = "'7.986"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.66%  "
$ws.Range("D49").Value = "'0.05222"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.50%  "
$ws.Range("D50").Value = "'0.4232"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.34%  "
$ws.Range("D51").Value = "'5.859"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.43%  "
